$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update contribution percentages (column D) and contribution descriptions (column E)
# for rows 9-12. Values are written in an order that reproduces the shared-string
# insertion order seen in the target workbook (E9, E10, E12, E11).

$ws.Range("D9").Value = 0.36
$ws.Range("E9").Value = "Organizou o grupo, implementou todo jogo e tabelas"

$ws.Range("D10").Value = 0.215
$ws.Range("E10").Value = "cuidou da página login"

$ws.Range("D12").Value = 0.225
$ws.Range("E12").Value = "cuidou da página registrar"

$ws.Range("D11").Value = 0.2
$ws.Range("E11").Value = "cuidou da página de alterar informações"

# Update the active selection to match the edited cell
$ws.Range("E13").Select()
